# Update (Removed Auto Arima)
# Updates the forecast figures on the "Forecast Comparison" sheet and the
# derived roll-up figures on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Forecast Comparison" - columns D (Amazon Mean Forecast),
# E (Amazon P70 Forecast), F (Amazon P80 Forecast), G (Amazon P90
# Forecast) for every data row, plus column C (Prophet Forecast) where
# it changed.
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

# row -> @(C, D, E, F, G)  (use $null to leave a column untouched)
$forecastRows = @{
    2  = @($null, 234, 277, 318, 382)
    3  = @(90,    198, 237, 276, 338)
    4  = @(91,    197, 235, 272, 330)
    5  = @(89,    200, 239, 278, 338)
    6  = @(74,    207, 249, 293, 361)
    7  = @(64,    197, 235, 275, 337)
    8  = @(61,    198, 239, 284, 354)
    9  = @(65,    198, 239, 283, 354)
    10 = @(72,    190, 228, 267, 329)
    11 = @(86,    188, 226, 268, 333)
    12 = @(98,    192, 232, 277, 346)
    13 = @(120,   192, 233, 283, 363)
    14 = @(147,   195, 237, 287, 367)
    15 = @($null, 190, 232, 284, 367)
    16 = @(127,   186, 227, 279, 362)
    17 = @(114,   185, 226, 276, 356)
}

foreach ($row in $forecastRows.Keys) {
    $vals = $forecastRows[$row]
    if ($null -ne $vals[0]) { $wsForecast.Cells.Item($row, 3).Value = $vals[0] }   # C
    if ($null -ne $vals[1]) { $wsForecast.Cells.Item($row, 4).Value = $vals[1] }   # D
    if ($null -ne $vals[2]) { $wsForecast.Cells.Item($row, 5).Value = $vals[2] }   # E
    if ($null -ne $vals[3]) { $wsForecast.Cells.Item($row, 6).Value = $vals[3] }   # F
    if ($null -ne $vals[4]) { $wsForecast.Cells.Item($row, 7).Value = $vals[4] }   # G
}

# ---------------------------------------------------------------------
# Sheet 2: "Summary" - roll-up metrics recomputed from the new forecast.
# These cells are stored as text, so a leading apostrophe is used to
# keep numeric-looking / date-looking strings from being re-interpreted
# as a number or a date serial.
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").Value  = "'1536"        # Total Forecast (16 Weeks)
$wsSummary.Range("B10").Value = "'618"         # Total Forecast (8 Weeks)
$wsSummary.Range("B11").Value = "'354"         # Total Forecast (4 Weeks)
$wsSummary.Range("B13").Value = "'2025-03-09"  # Max Forecast Week
$wsSummary.Range("B14").Value = "'61"          # Min Forecast
